$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("YDS")

$ws.Range("B2").Value = "MIE(47.33792324802306, 9.996989697741565, -25.92118911710068, 24.506680713814177)"
$ws.Range("C2").Value = "JSU(-1.2261633454245995, 0.9763695986829459, 2.7369937705124294, 3.5696905076586223)"
$ws.Range("D2").Value = "JSU(-0.8428691730978397, 1.0928559955026667, 0.6216691158531327, 2.556783385832321)"
$ws.Range("E2").Value = "NIG(0.7745491532141401, 0.5254706626336795, 5.6398847800278045, 5.193957139610973)"
